# 自动更新价格数据: insert a new row at the top of the data table (row 2)
# with the next day's date, pushing the existing history down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (the most recent data row).
# This shifts all existing data rows (2..75) down to (3..76).
$ws.Rows.Item(2).Insert()

# Populate the newly-inserted row 2 with the new date and the same
# commodity price figures as the rest of the (constant-valued) series.
# Force the date column to stay plain text (matches the existing rows,
# which are inline strings rather than numeric date serials).
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-02-03"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
